# Update Leave Card 4/12/2023 4:43 PM
# Adds 2023 year block (rows 78-107) to the leave card table on Sheet1,
# filling in VL usage for Mar 2023 and appending a new trailing table row
# so the table grows from A8:K130 to A8:K131.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow Table1 by one row -------------------------------------------
# Capture the formatting of the current (pre-insert) last table row (130)
# so we can re-apply it to the brand new last row (131) after the table
# auto-expands.
$ws.Range("A130:K130").Copy()
$null = $lo.ListRows.Add()
$ws.Range("A131:K131").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G131").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# The row that used to be the special "final" row (130) is now a normal
# interior row, so restore it to the standard row formatting (copied from
# row 129, an ordinary data row).
$ws.Range("A129:K129").Copy()
$ws.Range("A130:K130").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G130").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 2. Row 78: "2023" year divider --------------------------------------
# Match the formatting used by the other year-divider rows (e.g. row 10).
$ws.Range("A10").Copy()
$ws.Range("A78").PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A78").Formula = "'2023"

# --- 3. Monthly dates for the 2023 block (column A, rows 79-107) --------
$ws.Range("A79").Value2 = 44927
$ws.Range("A80").Value2 = 44958
$ws.Range("A81").Value2 = 44986
$ws.Range("A82").Value2 = 45017
$ws.Range("A83").Value2 = 45047
$ws.Range("A84").Value2 = 45078
$ws.Range("A85").Value2 = 45108
$ws.Range("A86").Value2 = 45139
$ws.Range("A87").Value2 = 45170
$ws.Range("A88").Value2 = 45200
$ws.Range("A89").Value2 = 45231
$ws.Range("A90").Value2 = 45261
$ws.Range("A91").Value2 = 45292
$ws.Range("A92").Value2 = 45323
$ws.Range("A93").Value2 = 45352
$ws.Range("A94").Value2 = 45383
$ws.Range("A95").Value2 = 45413
$ws.Range("A96").Value2 = 45444
$ws.Range("A97").Value2 = 45474
$ws.Range("A98").Value2 = 45505
$ws.Range("A99").Value2 = 45536
$ws.Range("A100").Value2 = 45566
$ws.Range("A101").Value2 = 45597
$ws.Range("A102").Value2 = 45627
$ws.Range("A103").Value2 = 45658
$ws.Range("A104").Value2 = 45689
$ws.Range("A105").Value2 = 45717
$ws.Range("A106").Value2 = 45748
$ws.Range("A107").Value2 = 45778

# --- 4. VL earned amounts for Jan-Mar 2023 (rows 79-81) ------------------
$ws.Range("C79").Value2 = 1.25
$ws.Range("C80").Value2 = 1.25
$ws.Range("C81").Value2 = 1.25

# --- 5. March 2023 leave usage entry (row 81) -----------------------------
$ws.Range("B81").Value2 = "VL(3-0-0)"
$ws.Range("D81").Value2 = 3
$ws.Range("K81").Value2 = "3/28,29,30/2023"

# --- 6. Update the saved view so the bottom pane shows the new rows ------
$ws.Activate()
$ws.Range("B82").Select()
